# Generate Report for Handoff
# Refresh the "Latest Handoff Date/Datetime" for the rows whose files were
# just re-handed-off: 3fec3e5f... (Handback transform failed) and the four
# "Ready for handoff" rows (41444385..., c496ae38..., e836c7da..., fab4bef2...).

$wb = $excel.ActiveWorkbook

# --- Overview sheet: column D = "Latest Handoff Date" ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("D4").Value  = "2016-03-19 04:31:37"
$wsOverview.Range("D7").Value  = "2016-03-19 04:31:37"
$wsOverview.Range("D8").Value  = "2016-03-19 04:31:37"
$wsOverview.Range("D9").Value  = "2016-03-19 04:31:37"
$wsOverview.Range("D10").Value = "2016-03-19 04:31:37"

# --- zh-cn sheet: column E = "Latest Handoff Datetime" ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value  = "2016-03-19 04:31:26"
$wsZhCn.Range("E7").Value  = "2016-03-19 04:31:26"
$wsZhCn.Range("E8").Value  = "2016-03-19 04:31:26"
$wsZhCn.Range("E9").Value  = "2016-03-19 04:31:26"
$wsZhCn.Range("E10").Value = "2016-03-19 04:31:26"

# --- de-de sheet: column E = "Latest Handoff Datetime" ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value  = "2016-03-19 04:31:37"
$wsDeDe.Range("E7").Value  = "2016-03-19 04:31:37"
$wsDeDe.Range("E8").Value  = "2016-03-19 04:31:37"
$wsDeDe.Range("E9").Value  = "2016-03-19 04:31:37"
$wsDeDe.Range("E10").Value = "2016-03-19 04:31:37"
